$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "last updated" timestamp string (row 1, column A)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 22:22"

# 2) Swap China / Brasil rows (rows 13 and 14) - country names
$ws.Range("A13").Value = "Brasil"
$ws.Range("A14").Value = "China"

# 3) Update numeric data for Estados Unidos (row 4)
$ws.Range("B4").Value = 1088415
$ws.Range("C4").Value = 24221
$ws.Range("D4").Value = 150768
$ws.Range("E4").Value = 874112
$ws.Range("F4").Value = 15226
$ws.Range("G4").Value = 1880
$ws.Range("H4").Value = 63535

# 4) Update numeric data for Alemania (row 9)
$ws.Range("B9").Value = 162530
$ws.Range("C9").Value = 991
$ws.Range("D9").Value = 123500
$ws.Range("E9").Value = 32458
$ws.Range("F9").Value = 2415
$ws.Range("G9").Value = 105
$ws.Range("H9").Value = 6572

# 5) Update numeric data for row 13 (now Brasil, with updated case counts)
$ws.Range("B13").Value = 85380
$ws.Range("C13").Value = 6019
$ws.Range("D13").Value = 34132
$ws.Range("E13").Value = 45347
$ws.Range("F13").Value = 8318
$ws.Range("G13").Value = 390
$ws.Range("H13").Value = 5901

# 6) Update numeric data for row 14 (now China, with prior data previously held by Brasil's old slot)
$ws.Range("B14").Value = 82862
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 77610
$ws.Range("E14").Value = 619
$ws.Range("F14").Value = 41
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 4633

# 7) Swap Namibia / San Vicente y las Granadinas rows (rows 191 and 192) - country names only
$ws.Range("A191").Value = "San Vicente y las Granadinas"
$ws.Range("A192").Value = "Namibia"

# 8) Update numeric data for Costa de Marfil (row 87)
$ws.Range("B87").Value = 1275
$ws.Range("C87").Value = 37
$ws.Range("D87").Value = 574
$ws.Range("E87").Value = 687
